$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 6 per the diff (Jogos_da_Semana_FlashScore_2024-11-13.xlsx)
$ws.Range("G6").Value = 1.6
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 6.1
$ws.Range("J6").Value = 2.1
$ws.Range("K6").Value = 2.07
$ws.Range("N6").Value = 6.85
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 2.67
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.62
$ws.Range("S6").Value = 1.42
$ws.Range("T6").Value = 2.47
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.65
$ws.Range("W6").Value = 5.5
$ws.Range("X6").Value = 6.7
$ws.Range("Y6").Value = 8
$ws.Range("Z6").Value = 11.5
$ws.Range("AA6").Value = 14
$ws.Range("AB6").Value = 32
$ws.Range("AC6").Value = 7.6
$ws.Range("AD6").Value = 6.7
$ws.Range("AE6").Value = 19.5
$ws.Range("AF6").Value = 120
$ws.Range("AG6").Value = 800
$ws.Range("AH6").Value = 13
$ws.Range("AI6").Value = 37
$ws.Range("AK6").Value = 150
$ws.Range("AL6").Value = 90
$ws.Range("AM6").Value = 90
$ws.Range("AO6").Value = 7.4
$ws.Range("AP6").Value = 17.5
$ws.Range("AQ6").Value = 24
$ws.Range("AS6").Value = 250
$ws.Range("AT6").Value = 2.45
$ws.Range("AU6").Value = 7.8
$ws.Range("AV6").Value = 80
$ws.Range("AW6").Value = 7.5
$ws.Range("AY6").Value = 45
$ws.Range("BA6").Value = 350
$ws.Range("BB6").Value = 450
